# "badania stat do hipotez" - build a small transposed legend table
# (variable name -> label) in columns AT:AU, rows 1-43, mirroring the
# header row (A1:AQ1) and the label row (A2:AQ2) of the questionnaire.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the two source rows (variable codes in row 1, human labels in row 2)
# across columns A:AQ (1..43) as a single 2 x 43 array.
$src = $ws.Range("A1:AQ2")
$vals = $src.Value2
$srcRows = $vals.GetLength(0)
$srcCols = $vals.GetLength(1)

# Build the transposed block: AT column gets row 1 (headers),
# AU column gets row 2 (labels) -> srcCols rows x srcRows columns.
$transposed = New-Object 'object[,]' $srcCols, $srcRows
for ($r = 1; $r -le $srcRows; $r++) {
    for ($c = 1; $c -le $srcCols; $c++) {
        $transposed[$c - 1, $r - 1] = $vals[$r, $c]
    }
}

$dst = $ws.Range("AT1:AU43")
$dst.Value2 = $transposed

# A2 (source of AU1) is formatted as Text; replicate that on the pasted cell.
$ws.Range("AU1").NumberFormat = "@"

# Manually widened column AT (46) after the paste, same as the source edit.
$ws.Columns.Item(46).ColumnWidth = 46 + 1/6

# Leave the selection where the editor ended up.
$ws.Range("AT27").Select() | Out-Null
